# The source data block (rows 186-258) shifts down by two rows, and two
# brand-new weekly price records are inserted at the top of that block
# (new rows 186 and 187). Use Excel's native row insert so every row
# below (188-260) is pushed down automatically, then populate the two
# freshly inserted rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 186 - this shifts the existing rows
# 186..258 down to 188..260 (and the dimension grows to A1:T260).
$ws.Rows.Item(186).Insert()
$ws.Rows.Item(186).Insert()

# --- New row 186 ---
$ws.Cells.Item(186, 1).Value = 4
$ws.Cells.Item(186, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value = "Los Lagos"
$ws.Cells.Item(186, 4).Value = 44726
$ws.Cells.Item(186, 5).Value = 10
$ws.Cells.Item(186, 6).Value = "Fruta"
$ws.Cells.Item(186, 7).Value = 100108
$ws.Cells.Item(186, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(186, 9).Value = 100108005
$ws.Cells.Item(186, 10).Value = "Piña"
$ws.Cells.Item(186, 11).Value = "Caramelo"
$ws.Cells.Item(186, 12).Value = "Especial"
$ws.Cells.Item(186, 13).Value = 100
$ws.Cells.Item(186, 14).Value = 18000
$ws.Cells.Item(186, 15).Value = 19000
$ws.Cells.Item(186, 16).Value = 18500
$ws.Cells.Item(186, 17).Value = "$/caja 10 unidades"
$ws.Cells.Item(186, 18).Value = "Ecuador"
$ws.Cells.Item(186, 19).Value = 1850
$ws.Cells.Item(186, 20).Value = 10

# --- New row 187 ---
$ws.Cells.Item(187, 1).Value = 4
$ws.Cells.Item(187, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(187, 3).Value = "Los Lagos"
$ws.Cells.Item(187, 4).Value = 44726
$ws.Cells.Item(187, 5).Value = 10
$ws.Cells.Item(187, 6).Value = "Fruta"
$ws.Cells.Item(187, 7).Value = 100108
$ws.Cells.Item(187, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(187, 9).Value = 100108005
$ws.Cells.Item(187, 10).Value = "Piña"
$ws.Cells.Item(187, 11).Value = "Caramelo"
$ws.Cells.Item(187, 12).Value = "Tercera"
$ws.Cells.Item(187, 13).Value = 160
$ws.Cells.Item(187, 14).Value = 20000
$ws.Cells.Item(187, 15).Value = 21000
$ws.Cells.Item(187, 16).Value = 20500
$ws.Cells.Item(187, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(187, 18).Value = "Ecuador"
$ws.Cells.Item(187, 19).Value = 1281
$ws.Cells.Item(187, 20).Value = 16

# Column D holds dates that must stay numeric (serial date values), not
# text - the rest of the column already uses this numeric date style.
$ws.Range("D186:D187").NumberFormat = $ws.Range("D188").NumberFormat
